$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 113.3752785
$ws.Range("H2").Value = 226.750557
$ws.Range("I2").Value = 0.06376184507388981
$ws.Range("J2").Value = 0.04522517177062929
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 7.930325
$ws.Range("N2").Value = 15.86065
$ws.Range("O2").Value = 0.1037243007717083
$ws.Range("P2").Value = 0.07315997304409949
$ws.Range("Q2").Value = 899.1028054705125
$ws.Range("R2").Value = 3596.41122188205
$ws.Range("S2").Value = 0.006613652796203216
$ws.Range("T2").Value = 0.003308672347654008

$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 113.3752785
$ws.Range("H3").Value = 226.750557
$ws.Range("I3").Value = 0.06376184507388981
$ws.Range("J3").Value = 0.04522517177062929
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 6.598259666666666
$ws.Range("N3").Value = 19.794779
$ws.Range("O3").Value = 0.08630161692429571
$ws.Range("P3").Value = 0.09130681895470279
$ws.Range("Q3").Value = 748.0795273236505
$ws.Range("R3").Value = 4488.477163941903
$ws.Range("S3").Value = 0.00550275032795313
$ws.Range("T3").Value = 0.004129366571056184

$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 113.3752785
$ws.Range("H4").Value = 226.750557
$ws.Range("I4").Value = 0.06376184507388981
$ws.Range("J4").Value = 0.04522517177062929
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 26.823721
$ws.Range("N4").Value = 80.471163
$ws.Range("O4").Value = 0.3508395563637543
$ws.Range("P4").Value = 0.371187064584827
$ws.Range("Q4").Value = 3041.146838781299
$ws.Range("R4").Value = 18246.88103268779
$ws.Range("S4").Value = 0.02237017743865793
$ws.Range("T4").Value = 0.01678699875488447

$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 113.3752785
$ws.Range("H5").Value = 226.750557
$ws.Range("I5").Value = 0.06376184507388981
$ws.Range("J5").Value = 0.04522517177062929
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 22.88541033333333
$ws.Range("N5").Value = 68.65623099999999
$ws.Range("O5").Value = 0.2993286132281626
$ws.Range("P5").Value = 0.3166886608852391
$ws.Range("Q5").Value = 2594.639770128444
$ws.Range("R5").Value = 15567.83862077067
$ws.Range("S5").Value = 0.01908574466283638
$ws.Range("T5").Value = 0.01432229908634551

$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 113.3752785
$ws.Range("H6").Value = 226.750557
$ws.Range("I6").Value = 0.06376184507388981
$ws.Range("J6").Value = 0.04522517177062929
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.575090333333333
$ws.Range("N6").Value = 22.725271
$ws.Range("O6").Value = 0.09907802619785785
$ws.Range("P6").Value = 0.1048242167742089
$ws.Range("Q6").Value = 858.8279762043245
$ws.Range("R6").Value = 5152.967857225947
$ws.Range("S6").Value = 0.006317397756654609
$ws.Range("T6").Value = 0.004740693209335276

$ws.Range("E7").Value = 2
$ws.Range("G7").Value = 113.3752785
$ws.Range("H7").Value = 226.750557
$ws.Range("I7").Value = 0.06376184507388981
$ws.Range("J7").Value = 0.04522517177062929
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 4.6429995
$ws.Range("N7").Value = 9.285999
$ws.Range("O7").Value = 0.06072788651422123
$ws.Range("P7").Value = 0.04283326575692263
$ws.Range("Q7").Value = 526.4013613878608
$ws.Range("R7").Value = 2105.605445551443
$ws.Range("S7").Value = 0.003872122091584536
$ws.Range("T7").Value = 0.00193714180135384

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1380.939473333333
$ws.Range("H8").Value = 4142.81842
$ws.Range("I8").Value = 0.7766353469649819
$ws.Range("J8").Value = 0.8262809897266405
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 7.930325
$ws.Range("N8").Value = 15.86065
$ws.Range("O8").Value = 0.1037243007717083
$ws.Range("P8").Value = 0.07315997304409949
$ws.Range("Q8").Value = 10951.29882886217
$ws.Range("R8").Value = 65707.792973173
$ws.Range("S8").Value = 0.08055595831853585
$ws.Range("T8").Value = 0.06045069493525287

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1380.939473333333
$ws.Range("H9").Value = 4142.81842
$ws.Range("I9").Value = 0.7766353469649819
$ws.Range("J9").Value = 0.8262809897266405
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 6.598259666666666
$ws.Range("N9").Value = 19.794779
$ws.Range("O9").Value = 0.08630161692429571
$ws.Range("P9").Value = 0.09130681895470279
$ws.Range("Q9").Value = 9111.797229003241
$ws.Range("R9").Value = 82006.17506102916
$ws.Range("S9").Value = 0.06702488620363935
$ws.Range("T9").Value = 0.075445088734683

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1380.939473333333
$ws.Range("H10").Value = 4142.81842
$ws.Range("I10").Value = 0.7766353469649819
$ws.Range("J10").Value = 0.8262809897266405
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 26.823721
$ws.Range("N10").Value = 80.471163
$ws.Range("O10").Value = 0.3508395563637543
$ws.Range("P10").Value = 0.371187064584827
$ws.Range("Q10").Value = 37041.93515058027
$ws.Range("R10").Value = 333377.4163552225
$ws.Range("S10").Value = 0.2724744005856047
$ws.Range("T10").Value = 0.3067048150988773

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1380.939473333333
$ws.Range("H11").Value = 4142.81842
$ws.Range("I11").Value = 0.7766353469649819
$ws.Range("J11").Value = 0.8262809897266405
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 22.88541033333333
$ws.Range("N11").Value = 68.65623099999999
$ws.Range("O11").Value = 0.2993286132281626
$ws.Range("P11").Value = 0.3166886608852391
$ws.Range("Q11").Value = 31603.36649273055
$ws.Range("R11").Value = 284430.298434575
$ws.Range("S11").Value = 0.2324691813910009
$ws.Range("T11").Value = 0.2616738201514598

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1380.939473333333
$ws.Range("H12").Value = 4142.81842
$ws.Range("I12").Value = 0.7766353469649819
$ws.Range("J12").Value = 0.8262809897266405
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 7.575090333333333
$ws.Range("N12").Value = 22.725271
$ws.Range("O12").Value = 0.09907802619785785
$ws.Range("P12").Value = 0.1048242167742089
$ws.Range("Q12").Value = 10460.74125536576
$ws.Range("R12").Value = 94146.67129829181
$ws.Range("S12").Value = 0.0769474972527789
$ws.Range("T12").Value = 0.08661425758351322

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1380.939473333333
$ws.Range("H13").Value = 4142.81842
$ws.Range("I13").Value = 0.7766353469649819
$ws.Range("J13").Value = 0.8262809897266405
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 4.6429995
$ws.Range("N13").Value = 9.285999
$ws.Range("O13").Value = 0.06072788651422123
$ws.Range("P13").Value = 0.04283326575692263
$ws.Range("Q13").Value = 6411.70128421693
$ws.Range("R13").Value = 38470.20770530157
$ws.Range("S13").Value = 0.04716342321342226
$ws.Range("T13").Value = 0.03539231322285425

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 1.722832
$ws.Range("H14").Value = 5.168496
$ws.Range("I14").Value = 0.000968914462885661
$ws.Range("J14").Value = 0.00103085135705228
$ws.Range("K14").Value = 2
$ws.Range("M14").Value = 7.930325
$ws.Range("N14").Value = 15.86065
$ws.Range("O14").Value = 0.1037243007717083
$ws.Range("P14").Value = 0.07315997304409949
$ws.Range("Q14").Value = 13.6626176804
$ws.Range("R14").Value = 81.9757060824
$ws.Range("S14").Value = 0.0001004999751704105
$ws.Range("T14").Value = 0.0000754170574944182

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 1.722832
$ws.Range("H15").Value = 5.168496
$ws.Range("I15").Value = 0.000968914462885661
$ws.Range("J15").Value = 0.00103085135705228
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 6.598259666666666
$ws.Range("N15").Value = 19.794779
$ws.Range("O15").Value = 0.08630161692429571
$ws.Range("P15").Value = 0.09130681895470279
$ws.Range("Q15").Value = 11.36769289804267
$ws.Range("R15").Value = 102.309236082384
$ws.Range("S15").Value = 0.00008361888480836805
$ws.Range("T15").Value = 0.00009412375822758224

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 1.722832
$ws.Range("H16").Value = 5.168496
$ws.Range("I16").Value = 0.000968914462885661
$ws.Range("J16").Value = 0.00103085135705228
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 26.823721
$ws.Range("N16").Value = 80.471163
$ws.Range("O16").Value = 0.3508395563637543
$ws.Range("P16").Value = 0.371187064584827
$ws.Range("Q16").Value = 46.21276489787201
$ws.Range("R16").Value = 415.914884080848
$ws.Range("S16").Value = 0.0003399335203132306
$ws.Range("T16").Value = 0.0003826386892475213

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 1.722832
$ws.Range("H17").Value = 5.168496
$ws.Range("I17").Value = 0.000968914462885661
$ws.Range("J17").Value = 0.00103085135705228
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 22.88541033333333
$ws.Range("N17").Value = 68.65623099999999
$ws.Range("O17").Value = 0.2993286132281626
$ws.Range("P17").Value = 0.3166886608852391
$ws.Range("Q17").Value = 39.42771725539733
$ws.Range("R17").Value = 354.8494552985759
$ws.Range("S17").Value = 0.0002900238225122749
$ws.Range("T17").Value = 0.0003264589358366181

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 1.722832
$ws.Range("H18").Value = 5.168496
$ws.Range("I18").Value = 0.000968914462885661
$ws.Range("J18").Value = 0.00103085135705228
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 7.575090333333333
$ws.Range("N18").Value = 22.725271
$ws.Range("O18").Value = 0.09907802619785785
$ws.Range("P18").Value = 0.1048242167742089
$ws.Range("Q18").Value = 13.05060802915733
$ws.Range("R18").Value = 117.455472262416
$ws.Range("S18").Value = 0.00009599813253726889
$ws.Range("T18").Value = 0.0001080581861136356

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 1.722832
$ws.Range("H19").Value = 5.168496
$ws.Range("I19").Value = 0.000968914462885661
$ws.Range("J19").Value = 0.00103085135705228
$ws.Range("K19").Value = 2
$ws.Range("M19").Value = 4.6429995
$ws.Range("N19").Value = 9.285999
$ws.Range("O19").Value = 0.06072788651422123
$ws.Range("P19").Value = 0.04283326575692263
$ws.Range("Q19").Value = 7.999108114584001
$ws.Range("R19").Value = 47.994648687504
$ws.Range("S19").Value = 0.00005884012754410804
$ws.Range("T19").Value = 0.00004415473013250466

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 4.12431
$ws.Range("H20").Value = 12.37293
$ws.Range("I20").Value = 0.002319496972672878
$ws.Range("J20").Value = 0.002467768511615926
$ws.Range("K20").Value = 2
$ws.Range("M20").Value = 7.930325
$ws.Range("N20").Value = 15.86065
$ws.Range("O20").Value = 0.1037243007717083
$ws.Range("P20").Value = 0.07315997304409949
$ws.Range("Q20").Value = 32.70711870075
$ws.Range("R20").Value = 196.2427122045
$ws.Range("S20").Value = 0.0002405882016325886
$ws.Range("T20").Value = 0.0001805418777888987

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 4.12431
$ws.Range("H21").Value = 12.37293
$ws.Range("I21").Value = 0.002319496972672878
$ws.Range("J21").Value = 0.002467768511615926
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 6.598259666666666
$ws.Range("N21").Value = 19.794779
$ws.Range("O21").Value = 0.08630161692429571
$ws.Range("P21").Value = 0.09130681895470279
$ws.Range("Q21").Value = 27.21326832583
$ws.Range("R21").Value = 244.91941493247
$ws.Range("S21").Value = 0.0002001763391926783
$ws.Range("T21").Value = 0.0002253240927122317

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 4.12431
$ws.Range("H22").Value = 12.37293
$ws.Range("I22").Value = 0.002319496972672878
$ws.Range("J22").Value = 0.002467768511615926
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 26.823721
$ws.Range("N22").Value = 80.471163
$ws.Range("O22").Value = 0.3508395563637543
$ws.Range("P22").Value = 0.371187064584827
$ws.Range("Q22").Value = 110.62934075751
$ws.Range("R22").Value = 995.6640668175901
$ws.Range("S22").Value = 0.0008137712888796238
$ws.Range("T22").Value = 0.0009160037499015833

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 4.12431
$ws.Range("H23").Value = 12.37293
$ws.Range("I23").Value = 0.002319496972672878
$ws.Range("J23").Value = 0.002467768511615926
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 22.88541033333333
$ws.Range("N23").Value = 68.65623099999999
$ws.Range("O23").Value = 0.2993286132281626
$ws.Range("P23").Value = 0.3166886608852391
$ws.Range("Q23").Value = 94.38652669186999
$ws.Range("R23").Value = 849.4787402268299
$ws.Range("S23").Value = 0.0006942918122170939
$ws.Range("T23").Value = 0.0007815143053184074

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 4.12431
$ws.Range("H24").Value = 12.37293
$ws.Range("I24").Value = 0.002319496972672878
$ws.Range("J24").Value = 0.002467768511615926
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 7.575090333333333
$ws.Range("N24").Value = 22.725271
$ws.Range("O24").Value = 0.09907802619785785
$ws.Range("P24").Value = 0.1048242167742089
$ws.Range("Q24").Value = 31.24202081267
$ws.Range("R24").Value = 281.17818731403
$ws.Range("S24").Value = 0.0002298111818243354
$ws.Range("T24").Value = 0.0002586819014101947

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 4.12431
$ws.Range("H25").Value = 12.37293
$ws.Range("I25").Value = 0.002319496972672878
$ws.Range("J25").Value = 0.002467768511615926
$ws.Range("K25").Value = 2
$ws.Range("M25").Value = 4.6429995
$ws.Range("N25").Value = 9.285999
$ws.Range("O25").Value = 0.06072788651422123
$ws.Range("P25").Value = 0.04283326575692263
$ws.Range("Q25").Value = 19.149169267845
$ws.Range("R25").Value = 114.89501560707
$ws.Range("S25").Value = 0.0001408581489265583
$ws.Range("T25").Value = 0.0001057025844846104

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 70.81572733333333
$ws.Range("H26").Value = 212.447182
$ws.Range("I26").Value = 0.03982650798977153
$ws.Range("J26").Value = 0.04237237793482528
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 7.930325
$ws.Range("N26").Value = 15.86065
$ws.Range("O26").Value = 0.1037243007717083
$ws.Range("P26").Value = 0.07315997304409949
$ws.Range("Q26").Value = 561.5917328647166
$ws.Range("R26").Value = 3369.5503971883
$ws.Range("S26").Value = 0.004130976693417908
$ws.Range("T26").Value = 0.003099962027526214

$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 70.81572733333333
$ws.Range("H27").Value = 212.447182
$ws.Range("I27").Value = 0.03982650798977153
$ws.Range("J27").Value = 0.04237237793482528
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 6.598259666666666
$ws.Range("N27").Value = 19.794779
$ws.Range("O27").Value = 0.08630161692429571
$ws.Range("P27").Value = 0.09130681895470279
$ws.Range("Q27").Value = 467.2605574291975
$ws.Range("R27").Value = 4205.345016862778
$ws.Range("S27").Value = 0.003437092035965665
$ws.Range("T27").Value = 0.003868887040775335

$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 70.81572733333333
$ws.Range("H28").Value = 212.447182
$ws.Range("I28").Value = 0.03982650798977153
$ws.Range("J28").Value = 0.04237237793482528
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 26.823721
$ws.Range("N28").Value = 80.471163
$ws.Range("O28").Value = 0.3508395563637543
$ws.Range("P28").Value = 0.371187064584827
$ws.Range("Q28").Value = 1899.541312401407
$ws.Range("R28").Value = 17095.87181161267
$ws.Range("S28").Value = 0.01397271439464896
$ws.Range("T28").Value = 0.01572807858510669

$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 70.81572733333333
$ws.Range("H29").Value = 212.447182
$ws.Range("I29").Value = 0.03982650798977153
$ws.Range("J29").Value = 0.04237237793482528
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 22.88541033333333
$ws.Range("N29").Value = 68.65623099999999
$ws.Range("O29").Value = 0.2993286132281626
$ws.Range("P29").Value = 0.3166886608852391
$ws.Range("Q29").Value = 1620.646978076782
$ws.Range("R29").Value = 14585.82280269104
$ws.Range("S29").Value = 0.01192121340629865
$ws.Range("T29").Value = 0.01341885162670307

$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 70.81572733333333
$ws.Range("H30").Value = 212.447182
$ws.Range("I30").Value = 0.03982650798977153
$ws.Range("J30").Value = 0.04237237793482528
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 7.575090333333333
$ws.Range("N30").Value = 22.725271
$ws.Range("O30").Value = 0.09907802619785785
$ws.Range("P30").Value = 0.1048242167742089
$ws.Range("Q30").Value = 536.4355315707023
$ws.Range("R30").Value = 4827.919784136322
$ws.Range("S30").Value = 0.003945931801979779
$ws.Range("T30").Value = 0.00444165132987883

$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 70.81572733333333
$ws.Range("H31").Value = 212.447182
$ws.Range("I31").Value = 0.03982650798977153
$ws.Range("J31").Value = 0.04237237793482528
$ws.Range("K31").Value = 2
$ws.Range("M31").Value = 4.6429995
$ws.Range("N31").Value = 9.285999
$ws.Range("O31").Value = 0.06072788651422123
$ws.Range("P31").Value = 0.04283326575692263
$ws.Range("Q31").Value = 328.797386600803
$ws.Range("R31").Value = 1972.784319604818
$ws.Range("S31").Value = 0.002418579657460571
$ws.Range("T31").Value = 0.001814947324835136

$ws.Range("E32").Value = 2
$ws.Range("G32").Value = 207.127739
$ws.Range("H32").Value = 414.255478
$ws.Range("I32").Value = 0.1164878885357982
$ws.Range("J32").Value = 0.08262284069923649
$ws.Range("K32").Value = 2
$ws.Range("M32").Value = 7.930325
$ws.Range("N32").Value = 15.86065
$ws.Range("O32").Value = 0.1037243007717083
$ws.Range("P32").Value = 0.07315997304409949
$ws.Range("Q32").Value = 1642.590286785175
$ws.Range("R32").Value = 6570.361147140699
$ws.Range("S32").Value = 0.01208262478674837
$ws.Range("T32").Value = 0.006044684798383068

$ws.Range("E33").Value = 2
$ws.Range("G33").Value = 207.127739
$ws.Range("H33").Value = 414.255478
$ws.Range("I33").Value = 0.1164878885357982
$ws.Range("J33").Value = 0.08262284069923649
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 6.598259666666666
$ws.Range("N33").Value = 19.794779
$ws.Range("O33").Value = 0.08630161692429571
$ws.Range("P33").Value = 0.09130681895470279
$ws.Range("Q33").Value = 1366.68260609156
$ws.Range("R33").Value = 8200.09563654936
$ws.Range("S33").Value = 0.01005309313273652
$ws.Range("T33").Value = 0.007544028757248435

$ws.Range("E34").Value = 2
$ws.Range("G34").Value = 207.127739
$ws.Range("H34").Value = 414.255478
$ws.Range("I34").Value = 0.1164878885357982
$ws.Range("J34").Value = 0.08262284069923649
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 26.823721
$ws.Range("N34").Value = 80.471163
$ws.Range("O34").Value = 0.3508395563637543
$ws.Range("P34").Value = 0.371187064584827
$ws.Range("Q34").Value = 5555.93668229682
$ws.Range("R34").Value = 33335.62009378091
$ws.Range("S34").Value = 0.04086855913564991
$ws.Range("T34").Value = 0.03066852970680937

$ws.Range("E35").Value = 2
$ws.Range("G35").Value = 207.127739
$ws.Range("H35").Value = 414.255478
$ws.Range("I35").Value = 0.1164878885357982
$ws.Range("J35").Value = 0.08262284069923649
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 22.88541033333333
$ws.Range("N35").Value = 68.65623099999999
$ws.Range("O35").Value = 0.2993286132281626
$ws.Range("P35").Value = 0.3166886608852391
$ws.Range("Q35").Value = 4740.203298430569
$ws.Range("R35").Value = 28441.21979058341
$ws.Range("S35").Value = 0.03486815813329727
$ws.Range("T35").Value = 0.02616571677957564

$ws.Range("E36").Value = 2
$ws.Range("G36").Value = 207.127739
$ws.Range("H36").Value = 414.255478
$ws.Range("I36").Value = 0.1164878885357982
$ws.Range("J36").Value = 0.08262284069923649
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 7.575090333333333
$ws.Range("N36").Value = 22.725271
$ws.Range("O36").Value = 0.09907802619785785
$ws.Range("P36").Value = 0.1048242167742089
$ws.Range("Q36").Value = 1569.01133346409
$ws.Range("R36").Value = 9414.068000784537
$ws.Range("S36").Value = 0.01154139007208296
$ws.Range("T36").Value = 0.008660874563957694

$ws.Range("E37").Value = 2
$ws.Range("G37").Value = 207.127739
$ws.Range("H37").Value = 414.255478
$ws.Range("I37").Value = 0.1164878885357982
$ws.Range("J37").Value = 0.08262284069923649
$ws.Range("K37").Value = 2
$ws.Range("M37").Value = 4.6429995
$ws.Range("N37").Value = 9.285999
$ws.Range("O37").Value = 0.06072788651422123
$ws.Range("P37").Value = 0.04283326575692263
$ws.Range("Q37").Value = 961.6939886131305
$ws.Range("R37").Value = 3846.775954452522
$ws.Range("S37").Value = 0.007074063275283208
$ws.Range("T37").Value = 0.00353900609326228
